# Tasks.xlsx edit: "fixed multilingual added standard template engine began
# with fixing soap service response interpretation"
#
# - Removed the stray "Partymarti" assignee string (replaced by "Kai"/"Riedo")
# - Added "obsolet" / "alle" status+assignee strings
# - Updated several rows' "Wer" (D) / "Status" (E) columns
# - Narrowed column B, cleared the scrolled/selected view state

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D ("Wer"/assignee) and column E ("Status") updates ---------

$ws.Range("E2").Value = "Done"

$ws.Range("E3").Value = "obsolet"
$ws.Range("E4").Value = "obsolet"

$ws.Range("D6").Value = "alle"

$ws.Range("D7").Value = "Kai"
$ws.Range("E7").Value = "done"

$ws.Range("D8").Value = "Kai"
$ws.Range("E8").Value = "done"

$ws.Range("D9").Value = "Alle"
$ws.Range("E9").Value = "done"

$ws.Range("D10").Value = "Kai"
$ws.Range("E10").Value = "done"

$ws.Range("E11").Value = "done"

$ws.Range("D12").Value = "Kai"
$ws.Range("E12").Value = "done"

$ws.Range("D13").Value = "Kai"
$ws.Range("E13").Value = "done"

$ws.Range("D14").Value = "Kai"
$ws.Range("E14").Value = "done"

$ws.Range("D15").Value = "Kai"
$ws.Range("E15").Value = "done"

$ws.Range("D16").Value = "Riedo"
$ws.Range("D17").Value = "Riedo"
$ws.Range("D18").Value = "Riedo"
$ws.Range("D19").Value = "Riedo"
$ws.Range("D20").Value = "Riedo"
$ws.Range("D21").Value = "Riedo"

$ws.Range("D22").Value = "Alle"

$ws.Range("D23").Value = "Riedo"
$ws.Range("D24").Value = "Riedo"
$ws.Range("D25").Value = "Riedo"

$ws.Range("D26").Value = "Kai"

$ws.Range("D27").Value = "Riedo"

$ws.Range("D28").Value = "Alle"

$ws.Range("D29").Value = "Kai"
$ws.Range("E29").Value = "done"

$ws.Range("D30").Value = "Kai"
$ws.Range("E30").Value = "done"

$ws.Range("D31").Value = "Kai"
$ws.Range("E31").Value = "done"

$ws.Range("D32").Value = "Kai"
$ws.Range("E32").Value = "begonnen"

$ws.Range("D33").Value = "Kai"
$ws.Range("E33").Value = "begonnen"

$ws.Range("D34").Value = "Riedo"

# --- View state: narrower "Wer" column, cursor parked at E5, no frozen ---
# --- scroll position (topLeftCell) --------------------------------------

$ws.Columns.Item(2).ColumnWidth = 9.5
$ws.Range("E5").Select()
